# Update countries & provincias Spain
#
# The sheet "Pais" lists countries ranked by Casos totales (column B,
# descending). This applies the day's data refresh:
#   - Most rows just get new totals (Casos totales/Nuevos casos/Casos
#     activos/Recuperados/Casos criticos/Muertes hoy/Muertes) in B:H.
#   - "Nepal" newly overtakes "Bahamas" in the ranking, so its figures
#     now occupy the row Bahamas used to have, pushing Bahamas..San
#     Martin (Parte Holandesa) down by one rank each.
#   - "Republica Dominicana" grows past "Noruega" and "Chequia", so it
#     now occupies the row Noruega used to have, and Noruega/Chequia
#     drop one rank each.
# Column A (country name) + B:H (stats) are written together per row so
# every row stays internally consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=4;   Values=@('Estados Unidos',               1171350, 10576,  174764, 928498, 16366, 644, 68088)}
    @{Row=18;  Values=@('India',                         42505,  2806,   11775,  29339,  0,     68,  1391)}
    @{Row=26;  Values=@('Irlanda',                        21506,  330,    13386,  6817,   99,    17,  1303)}
    @{Row=28;  Values=@('Chile',                          19663,  1228,   10041,  9362,   408,   13,  260)}
    @{Row=43;  Values=@('Banglades',                       9455,  665,    1063,   8215,   1,     2,   177)}
    @{Row=45;  Values=@('Republica Dominicana',            7954,  376,    1606,   6015,   144,   7,   333)}
    @{Row=46;  Values=@('Noruega',                         7809,  0,      32,     7566,   37,    0,   211)}
    @{Row=47;  Values=@('Chequia',                         7764,  9,      3584,   3935,   62,    0,   245)}
    @{Row=70;  Values=@('Irak',                            2296,  77,     1490,   709,    0,     2,   97)}
    @{Row=160; Values=@('Nepal',                             84,  25,     16,     68,     0,     0,   0)}
    @{Row=161; Values=@('Bahamas',                           83,  0,      24,     48,     1,     0,   11)}
    @{Row=162; Values=@('Guyana',                            82,  0,      22,     51,     2,     0,   9)}
    @{Row=163; Values=@('Liechtenstein',                     82,  0,      55,     26,     0,     0,   1)}
    @{Row=164; Values=@('Barbados',                          81,  0,      44,     30,     4,     0,   7)}
    @{Row=165; Values=@('Mozambique',                        79,  0,      18,     61,     0,     0,   0)}
    @{Row=166; Values=@('San Martin (Parte Holandesa)',      76,  0,      44,     19,     7,     0,   13)}
)

foreach ($update in $updates) {
    $r = $update.Row
    $values = $update.Values
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $values[$col - 1]
    }
}
